$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Backfill two new daily rows (8 April and 9 April 2020) from ESR data.
# Columns: A date, B confirmed, C totalConfirmed, D probable, E totalProbable,
# F cases, G totalCases, H recovered, I totalRecovered, J inHospitalNow,
# K totalBeenInHospital (not populated), L inIcu, M deaths, N totalDeaths,
# O overseas, P contact, Q investigating, R community, S established, T tag

$newRows = @(
    @{ Row = 42; A = 43929; B = 26; C = 969; D = 24; E = 241; F = 50; G = 1210; H = 41; I = 282; J = 12; L = 4; M = 0; N = 1; O = 496; P = 520; Q = 169; R = 24; S = 1210; T = "Manual" },
    @{ Row = 43; A = 43930; B = 23; C = 992; D = 6;  E = 247; F = 29; G = 1239; H = 35; I = 317; J = 14; L = 4; M = 0; N = 1; O = 508; P = 545; Q = 161; R = 25; S = 1239; T = "Manual" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 1).NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    # Column K (totalBeenInHospital) intentionally left blank
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
